$d = $word.ActiveDocument

# Target colors (Word BGR-packed RGB integers: r + g*256 + b*65536)
$red  = 3355647   # RGB(0xFF,0x33,0x33) -> "FF3333"
$blue = 16750899  # RGB(0x33,0x99,0xFF) -> "3399FF"

# The "Function 1".."Function 4" block runs from paragraph 14 through
# paragraph 43 (1-based Paragraphs collection). Recolor the *text* of each
# paragraph (not the trailing paragraph-mark) from #333333 to #FF3333,
# except the "Description: ..." paragraph right after "Function 1", which
# becomes #3399FF.
$firstPara = 14
$lastPara  = 43
$blueParaText = "Description:  Function takes a hostname, determines the IP address(es) for the host and pings each IP address to determine if it is online.  Return output that shows results of ping."
$specialParaText = "-Count:  Optional Number of times to ping the device"

for ($i = $firstPara; $i -le $lastPara; $i++) {
    $p = $d.Paragraphs.Item($i)
    # Paragraph.Range.Text includes the trailing paragraph-mark (CR)
    # character; strip it before comparing against plain text.
    $txt = $p.Range.Text.TrimEnd([char]13)

    if ($txt -eq $specialParaText) {
        # This paragraph also gets its paragraph-mark (end-of-paragraph)
        # character formatting recolored, in addition to its text.
        $p.Range.Font.Color = $red
        continue
    }

    # Text-only range: excludes the trailing paragraph mark so the
    # paragraph-mark run properties (w:pPr/w:rPr) are left untouched.
    $textRange = $d.Range($p.Range.Start, $p.Range.End - 1)

    if ($txt -eq $blueParaText) {
        $textRange.Font.Color = $blue
    } else {
        $textRange.Font.Color = $red
    }
}

Write-Host "done"
